$wb = $excel.ActiveWorkbook

$totalSheet = $wb.Worksheets.Item(1)
$q2Sheet = $wb.Worksheets.Item(2)

# ---------------------------------------------------------------------------
# 1. Insert the new "2022-Q3" sheet right after "总计" (i.e. before the
#    existing "2022-Q2" sheet), mirroring the new tab order in the diff.
# ---------------------------------------------------------------------------
$newSheet = $wb.Worksheets.Add($q2Sheet)
$newSheet.Name = "2022-Q3"

# Re-use the bold/bordered header style (style index 2 in the original file,
# carried by 总计!B1) for the new sheet's header row and index column, so we
# don't introduce a bunch of near-duplicate styles.
$totalSheet.Range("B1").Copy($newSheet.Range("B1:H1"))
$totalSheet.Range("A2").Copy($newSheet.Range("A2:A4"))

# Header row
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# Index column (numeric, 0-based row counter)
$newSheet.Range("A2").Value = 0
$newSheet.Range("A3").Value = 1
$newSheet.Range("A4").Value = 2

# Fund code / name columns - B is numeric-looking text, C is plain text.
$newSheet.Range("B2").NumberFormat = "@"
$newSheet.Range("B3").NumberFormat = "@"
$newSheet.Range("B4").NumberFormat = "@"
$newSheet.Range("B2").Value = "519991"
$newSheet.Range("C2").Value = "长信双利优选混合A"
$newSheet.Range("B3").Value = "168301"
$newSheet.Range("C3").Value = "东海祥龙灵活配置混合（LOF）"
$newSheet.Range("B4").Value = "006396"
$newSheet.Range("C4").Value = "长信双利优选混合E"

# Numeric-looking values stored as text (D, E, F, G columns) - force text so
# they keep their original formatting (e.g. "0.00", "91.61") instead of being
# parsed as numbers.
$newSheet.Range("D2:G3").NumberFormat = "@"
$newSheet.Range("D4:F4").NumberFormat = "@"

$newSheet.Range("D2").Value = "1.06"
$newSheet.Range("E2").Value = "91.61"
$newSheet.Range("F2").Value = "5.98"
$newSheet.Range("G2").Value = "0.0634"
$newSheet.Range("H2").Value = 2

$newSheet.Range("D3").Value = "0.14"
$newSheet.Range("E3").Value = "87.75"
$newSheet.Range("F3").Value = "2.56"
$newSheet.Range("G3").Value = "0.0036"
$newSheet.Range("H3").Value = 9

$newSheet.Range("D4").Value = "0.00"
$newSheet.Range("E4").Value = "91.61"
$newSheet.Range("F4").Value = "5.98"
$newSheet.Range("G4").Value = 0
$newSheet.Range("H4").Value = 2

# ---------------------------------------------------------------------------
# 2. Update the "总计" summary sheet: the new quarter is inserted as row 2
#    and every later row shifts down by one (their B/C/D content is
#    unchanged, only the running index in column A is renumbered).
# ---------------------------------------------------------------------------

# Make room for the extra row by copying the last row's format down first
# (keeps column A's bold/bordered style on the brand-new row 9), then fill in
# the values top to bottom.
$totalSheet.Range("A8").Copy($totalSheet.Range("A9"))

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q3"
$totalSheet.Range("C2").Value = 3
$totalSheet.Range("D2").Value = 0.07000000000000001

$totalSheet.Range("A3").Value = 1
$totalSheet.Range("B3").Value = "2022-Q2"
$totalSheet.Range("C3").Value = 6
$totalSheet.Range("D3").Value = 0.08

$totalSheet.Range("A4").Value = 2
$totalSheet.Range("B4").Value = "2022-Q1"
$totalSheet.Range("C4").Value = 4
$totalSheet.Range("D4").Value = 0.3

$totalSheet.Range("A5").Value = 3
$totalSheet.Range("B5").Value = "2021-Q4"
$totalSheet.Range("C5").Value = 2
$totalSheet.Range("D5").Value = 0.17

$totalSheet.Range("A6").Value = 4
$totalSheet.Range("B6").Value = "2021-Q3"
$totalSheet.Range("C6").Value = 6
$totalSheet.Range("D6").Value = 0.71

$totalSheet.Range("A7").Value = 5
$totalSheet.Range("B7").Value = "2021-Q2"
$totalSheet.Range("C7").Value = 1
$totalSheet.Range("D7").Value = 0.23

$totalSheet.Range("A8").Value = 6
$totalSheet.Range("B8").Value = "2021-Q1"
$totalSheet.Range("C8").Value = 3
$totalSheet.Range("D8").Value = 0.28

$totalSheet.Range("A9").Value = 7
$totalSheet.Range("B9").Value = "2020-Q4"
$totalSheet.Range("C9").Value = 4
$totalSheet.Range("D9").Value = 0.49
